# Apply UC007 spreadsheet update (v1.2.1 -> v1.2.3)
# 1) Update the "Expected Results" text for the first step (list pending daily allowances)
#    in every test case (TC1..TC5, which all share this text) to mention ordering by
#    diária number ascending.
# 2) Swap the second-step content between TC2 and TC4:
#    - TC2's second step becomes "filter by user" (previously TC4's second step)
#    - TC4's second step becomes "assign/unassign" (previously TC2's second step)
#    TC3 stays unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newListText = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo numero de diarias em ordem crescente."

# Rows whose column D holds the "list of dailies" expected-result text (one per test case block:
# TC1, TC2, TC3, TC4, TC5 all share this same text).
$listRows = @(10, 19, 27, 35, 43)
foreach ($r in $listRows) {
    $ws.Cells.Item($r, 4).Value = $newListText
}

# Text used for the "assign/unassign" step (previously TC2's second step)
$assignStep = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$assignResult = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# Text used for the "filter by user" step (previously TC4's second step)
$filterStep = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$filterResult = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."

# TC2 second step (row 20): was "assign/unassign" -> becomes "filter by user"
$ws.Cells.Item(20, 2).Value = $filterStep
$ws.Cells.Item(20, 4).Value = $filterResult

# TC4 second step (row 36): was "filter by user" -> becomes "assign/unassign"
$ws.Cells.Item(36, 2).Value = $assignStep
$ws.Cells.Item(36, 4).Value = $assignResult
